# Deploy IG build 2e04aec7a6e183fa6df6ea6f057cb991a5ed746d
# Sets the "Experimental" and "Case Sensitive" metadata values on the
# Metadata worksheet, which were previously blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B7").Value = "false"
$ws.Range("B14").Value = "true"
